# PlayerUpgradeDB_Sheet.xlsx - minor fixes
# - Reset the "cost"-looking column F values (rows 2-5) from 1234 to 0
# - Move the active selection from Q10 to L9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Zero out column F (rows 2-5)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0

# Update the saved selection/active cell on the sheet
$ws.Activate()
[void]$ws.Range("L9").Select()
